$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
# (Sheet1 is also the workbook's ActiveSheet in this file.)

$ws.Range("A28").Value = 16
$ws.Range("B28").Value = "Unison"
$ws.Range("C28").Value = "For syncing folders."

# Reflect the post-edit UI state (cursor moved on to the next row after
# typing the final cell of the new entry).
[void]$ws.Range("C29").Select()
